# ------------------------------------------------------------------
# 5-56.xlsx : build "Statement of Cash Flow" sheet content
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Statement of Cash Flow"

# ------------------------------------------------------------------
# Shared-string text constants (indices mirror the target sharedStrings.xml)
# ------------------------------------------------------------------
$s0  = 'Particulars'
$s1  = 'Amount'
$s2  = 'Component'
$s3  = 'Total'
$s4  = 'Name : Jacinta Manufacturing Co. (JMC)'
$s5  = 'Cash flow by Operations'
$s6  = '           Cash Sales'
$s7  = '           Received from Outstanding A/c'
$s8  = '           Operation Expenses:'
$s9  = '                   Purchase of goods on cash'
$s10 = '                   Salaries expense'
$s11 = '                   Other Operational expenses'
$s12 = '           Income Tax Applied'
$s13 = 'Total from Operations'
$s14 = 'Cash flow by Investments'
$s15 = '            Purchase of welding machines'
$s16 = '            Sales of old stock'
$s17 = 'Total from Investments'
$s18 = 'Cash flow from Finance'
$s19 = '            Issued a Long Term Debt'
$s20 = 'Total from Finance'
$s21 = 'Net cash flow from activities'
$s22 = 'Cash & Cash Equivalents on 31 Dec, 2010'
$s23 = 'Cash & Cash Equivalents on 31 Dec, 2011'
$s24 = '           Total Operational Expenses '
$s25 = 'STATEMENT OF CASH FLOW AS OF 31 DEC, 2011 (in thousands of $)'

# ------------------------------------------------------------------
# Column widths (best achievable approximation - the runtime snaps
# column widths to its own internal character grid)
# ------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth  = 24.75
$ws.Columns.Item(10).ColumnWidth = 10.251
$ws.Columns.Item(11).ColumnWidth = 1.251
$ws.Columns.Item(12).ColumnWidth = 1.75
$ws.Columns.Item(13).ColumnWidth = 12.084
$ws.Columns.Item(14).ColumnWidth = 11.75

# ------------------------------------------------------------------
# Title
# ------------------------------------------------------------------
$ws.Range("K5").Value = $s25
$ws.Range("K5").Font.Bold = $true
$ws.Range("K5").HorizontalAlignment = -4108   # xlCenter

# ------------------------------------------------------------------
# Company name
# ------------------------------------------------------------------
$ws.Range("I6").Value = $s4

# ------------------------------------------------------------------
# Amount (merged M7:N7) header
# ------------------------------------------------------------------
$ws.Range("M7:N7").Merge()
$ws.Range("M7").Value = $s1
$ws.Range("M7:N7").Font.Bold = $true
$ws.Range("M7:N7").HorizontalAlignment = -4108   # xlCenter

# ------------------------------------------------------------------
# Column headers
# ------------------------------------------------------------------
$ws.Range("I8").Value = $s0
$ws.Range("M8").Value = $s2
$ws.Range("N8").Value = $s3
$ws.Range("I8").Font.Bold = $true
$ws.Range("I8").HorizontalAlignment = -4108   # xlCenter
$ws.Range("M8").Font.Bold = $true
$ws.Range("M8").HorizontalAlignment = -4108   # xlCenter
$ws.Range("N8").Font.Bold = $true
$ws.Range("N8").HorizontalAlignment = -4108   # xlCenter

# ------------------------------------------------------------------
# Operations section
# ------------------------------------------------------------------
$ws.Range("I10").Value = $s5
$ws.Range("I10").Font.Bold = $true

$ws.Range("I11").Value = $s6
$ws.Range("N11").Value = 490

$ws.Range("I12").Value = $s7
$ws.Range("N12").Value = 15

$ws.Range("I13").Value = $s8

$ws.Range("I14").Value = $s9
$ws.Range("M14").Value = -300

$ws.Range("I15").Value = $s10
$ws.Range("M15").Value = -82

$ws.Range("I16").Value = $s11
$ws.Range("M16").Value = -15

$ws.Range("I17").Value = $s24
$ws.Range("N17").Formula = "=SUM(M14:M16)"

$ws.Range("I18").Value = $s12
$ws.Range("N18").Value = -8

$ws.Range("I20").Value = $s13
$ws.Range("I20").Font.Bold = $true
$ws.Range("N20").Formula = "=SUM(N11:N18)"
$ws.Range("N20").Font.Bold = $true

# ------------------------------------------------------------------
# Investments section
# ------------------------------------------------------------------
$ws.Range("I22").Value = $s14
$ws.Range("I22").Font.Bold = $true

$ws.Range("I23").Value = $s15
$ws.Range("M23").Value = -125

$ws.Range("I24").Value = $s16
$ws.Range("M24").Value = 5

$ws.Range("I26").Value = $s17
$ws.Range("I26").Font.Bold = $true
$ws.Range("N26").Formula = "=SUM(M23:M24)"
$ws.Range("N26").Font.Bold = $true

# ------------------------------------------------------------------
# Finance section
# ------------------------------------------------------------------
$ws.Range("I28").Value = $s18
$ws.Range("I28").Font.Bold = $true

$ws.Range("I29").Value = $s19
$ws.Range("M29").Value = 100

$ws.Range("I31").Value = $s20
$ws.Range("I31").Font.Bold = $true
$ws.Range("N31").Formula = "=SUM(M29)"
$ws.Range("N31").Font.Bold = $true

# ------------------------------------------------------------------
# Totals
# ------------------------------------------------------------------
$ws.Range("I33").Value = $s21
$ws.Range("I33").Font.Bold = $true
$ws.Range("N33").Formula = "=SUM(N20, N26, N31)"
$ws.Range("N33").Font.Bold = $true

$ws.Range("I34").Value = $s22
$ws.Range("N34").Value = 45

$ws.Range("I36").Value = $s23
$ws.Range("I36").Font.Bold = $true
$ws.Range("N36").Formula = "=SUM(N33:N34)"
$ws.Range("N36").Font.Bold = $true

# ------------------------------------------------------------------
# Selection, matching the target sheetView
# ------------------------------------------------------------------
$ws.Range("N22").Select()
